# Apply the "Gaussian Quadrature Scheme" export update to the averaged-intensities worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet (drives <sheets><sheet name="...">).
$ws.Name = "alpha4F"

# 2. Correct a tiny floating point rounding difference in I13.
$ws.Range("I13").Value = 0.9937878340704579

# 3. Append a new row of averaged intensity data (row 16) for the
#    "HexGrid-60degTilt5degRes" scheme (index 14), matching the formatting
#    used by the other index cells in column A (row 15 as a template).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9975163270514803
$ws.Range("D16").Value = 0.9929265524063458
$ws.Range("E16").Value = 0.9988235294117647
$ws.Range("F16").Value = 0.9975163270514803
$ws.Range("G16").Value = 0.9983006600036716
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 0.9976470588235294
$ws.Range("J16").Value = 0.9929265524063458
$ws.Range("K16").Value = 0.9958750409090552
$ws.Range("L16").Value = 0.9966956839802676
$ws.Range("M16").Value = 0.9975356879494653
